$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = 0.66045457039765831
$ws.Range("D2").Value = 0.1796404977586899
$ws.Range("E2").Value = 0.98377283629348966

# Row 3 updates
$ws.Range("B3").Value = 0.13420939667702694
$ws.Range("C3").Value = 1.8271309109788068
$ws.Range("D3").Value = 0.1508268078743738
$ws.Range("E3").Value = 1.2538360540595634

# Update selection to match new range
$ws.Range("B1:E3").Select()
